$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "button_closeAction_class"
$ws.Range("B1").Value = "button_closeAction_class_1"
$ws.Range("C1").Value = "div_dynamicObject_class"
$ws.Range("D1").Value = "div_dynamicObject_class_1"
